$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Subgroups")

$ws.Range("A22").Value = '{''Hobby'': ''1'', ''Student'': ''1''}'
$ws.Range("A23").Value = '{''FormalEducation'': ''1'', ''Hobby'': ''1''}'
$ws.Range("A25").Value = '{''Hobby'': ''1'', ''DevType'': ''2''}'
$ws.Range("A27").Value = '{''Hobby'': ''1'', ''SexualOrientation'': ''1''}'
$ws.Range("A29").Value = '{''Dependents'': ''2'', ''Hobby'': ''1''}'
$ws.Range("A31").Value = '{''Age'': ''3'', ''Hobby'': ''1''}'
$ws.Range("A33").Value = '{''Hobby'': ''1'', ''GDP'': ''1''}'
$ws.Range("A37").Value = '{''Country'': ''2'', ''SexualOrientation'': ''1''}'
$ws.Range("A40").Value = '{''Country'': ''2'', ''GDP'': ''2''}'
$ws.Range("A42").Value = '{''FormalEducation'': ''1'', ''Student'': ''1''}'
$ws.Range("A43").Value = '{''UndergradMajor'': ''2'', ''Student'': ''1''}'
$ws.Range("A44").Value = '{''Student'': ''1'', ''DevType'': ''2''}'
$ws.Range("A45").Value = '{''Exercise'': ''3'', ''Student'': ''1''}'
$ws.Range("A48").Value = '{''Dependents'': ''2'', ''Student'': ''1''}'
$ws.Range("A50").Value = '{''Age'': ''3'', ''Student'': ''1''}'
$ws.Range("A51").Value = '{''HDI'': ''1'', ''Student'': ''1''}'
$ws.Range("A52").Value = '{''GDP'': ''1'', ''Student'': ''1''}'
$ws.Range("A54").Value = '{''GINI'': ''1'', ''Student'': ''1''}'
$ws.Range("A55").Value = '{''GINI'': ''2'', ''Student'': ''1''}'
$ws.Range("A56").Value = '{''FormalEducation'': ''1'', ''UndergradMajor'': ''2''}'
$ws.Range("A57").Value = '{''FormalEducation'': ''1'', ''DevType'': ''2''}'
$ws.Range("A58").Value = '{''FormalEducation'': ''1'', ''SexualOrientation'': ''1''}'
$ws.Range("A59").Value = '{''FormalEducation'': ''1'', ''RaceEthnicity'': ''1''}'
$ws.Range("A60").Value = '{''FormalEducation'': ''1'', ''Dependents'': ''2''}'
$ws.Range("A61").Value = '{''FormalEducation'': ''1'', ''HDI'': ''1''}'
$ws.Range("A63").Value = '{''UndergradMajor'': ''2'', ''DevType'': ''2''}'
$ws.Range("A64").Value = '{''UndergradMajor'': ''2'', ''SexualOrientation'': ''1''}'
$ws.Range("A65").Value = '{''UndergradMajor'': ''2'', ''RaceEthnicity'': ''1''}'
$ws.Range("A67").Value = '{''Age'': ''3'', ''UndergradMajor'': ''2''}'
$ws.Range("A69").Value = '{''UndergradMajor'': ''2'', ''GDP'': ''1''}'
$ws.Range("A73").Value = '{''Dependents'': ''2'', ''DevType'': ''2''}'
$ws.Range("A74").Value = '{''Age'': ''3'', ''DevType'': ''2''}'
$ws.Range("A75").Value = '{''HDI'': ''1'', ''DevType'': ''2''}'
$ws.Range("A76").Value = '{''GINI'': ''2'', ''DevType'': ''2''}'
$ws.Range("A77").Value = '{''Exercise'': ''3'', ''SexualOrientation'': ''1''}'
$ws.Range("A79").Value = '{''SexualOrientation'': ''1'', ''Dependents'': ''1''}'
$ws.Range("A80").Value = '{''Dependents'': ''2'', ''SexualOrientation'': ''1''}'
$ws.Range("A81").Value = '{''Continent'': ''1'', ''SexualOrientation'': ''1''}'
$ws.Range("A82").Value = '{''Age'': ''3'', ''SexualOrientation'': ''1''}'
$ws.Range("A83").Value = '{''HDI'': ''1'', ''SexualOrientation'': ''1''}'
$ws.Range("A86").Value = '{''GINI'': ''2'', ''SexualOrientation'': ''1''}'
$ws.Range("A89").Value = '{''Age'': ''3'', ''RaceEthnicity'': ''1''}'
$ws.Range("A92").Value = '{''GINI'': ''1'', ''RaceEthnicity'': ''1''}'
$ws.Range("A93").Value = '{''GINI'': ''2'', ''RaceEthnicity'': ''1''}'
$ws.Range("A94").Value = '{''Dependents'': ''2'', ''Continent'': ''1''}'
$ws.Range("A95").Value = '{''Age'': ''3'', ''Dependents'': ''2''}'
$ws.Range("A96").Value = '{''Dependents'': ''2'', ''HDI'': ''1''}'
$ws.Range("A97").Value = '{''Dependents'': ''2'', ''GDP'': ''1''}'
$ws.Range("A99").Value = '{''HDI'': ''1'', ''Continent'': ''1''}'
$ws.Range("A102").Value = '{''GINI'': ''2'', ''Continent'': ''2''}'
$ws.Range("A103").Value = '{''Age'': ''3'', ''HDI'': ''1''}'
$ws.Range("A104").Value = '{''HDI'': ''1'', ''GDP'': ''1''}'
$ws.Range("A105").Value = '{''HDI'': ''1'', ''GDP'': ''2''}'
$ws.Range("A108").Value = '{''GINI'': ''1'', ''GDP'': ''1''}'
$ws.Range("A110").Value = '{''FormalEducation'': ''1'', ''Hobby'': ''1'', ''Student'': ''1''}'
$ws.Range("A111").Value = '{''UndergradMajor'': ''2'', ''Hobby'': ''1'', ''Student'': ''1''}'
$ws.Range("A112").Value = '{''Hobby'': ''1'', ''Student'': ''1'', ''DevType'': ''2''}'
$ws.Range("A113").Value = '{''Hobby'': ''1'', ''Student'': ''1'', ''SexualOrientation'': ''1''}'
$ws.Range("A114").Value = '{''RaceEthnicity'': ''1'', ''Hobby'': ''1'', ''Student'': ''1''}'
$ws.Range("A115").Value = '{''Dependents'': ''2'', ''Hobby'': ''1'', ''Student'': ''1''}'
$ws.Range("A116").Value = '{''Hobby'': ''1'', ''Student'': ''1'', ''Continent'': ''1''}'
$ws.Range("A117").Value = '{''Age'': ''3'', ''Hobby'': ''1'', ''Student'': ''1''}'
$ws.Range("A118").Value = '{''Hobby'': ''1'', ''Student'': ''1'', ''HDI'': ''1''}'
$ws.Range("A119").Value = '{''Hobby'': ''1'', ''Student'': ''1'', ''GDP'': ''1''}'
$ws.Range("A120").Value = '{''GINI'': ''2'', ''Hobby'': ''1'', ''Student'': ''1''}'
$ws.Range("A121").Value = '{''FormalEducation'': ''1'', ''Hobby'': ''1'', ''SexualOrientation'': ''1''}'
$ws.Range("A122").Value = '{''FormalEducation'': ''1'', ''Hobby'': ''1'', ''RaceEthnicity'': ''1''}'
$ws.Range("A123").Value = '{''FormalEducation'': ''1'', ''Hobby'': ''1'', ''HDI'': ''1''}'
$ws.Range("A124").Value = '{''UndergradMajor'': ''2'', ''Hobby'': ''1'', ''DevType'': ''2''}'
$ws.Range("A125").Value = '{''UndergradMajor'': ''2'', ''Hobby'': ''1'', ''SexualOrientation'': ''1''}'
$ws.Range("A126").Value = '{''UndergradMajor'': ''2'', ''Hobby'': ''1'', ''RaceEthnicity'': ''1''}'
$ws.Range("A129").Value = '{''Hobby'': ''1'', ''DevType'': ''2'', ''SexualOrientation'': ''1''}'
$ws.Range("A130").Value = '{''RaceEthnicity'': ''1'', ''Hobby'': ''1'', ''DevType'': ''2''}'
$ws.Range("A131").Value = '{''Dependents'': ''2'', ''Hobby'': ''1'', ''DevType'': ''2''}'
$ws.Range("A132").Value = '{''Hobby'': ''1'', ''DevType'': ''2'', ''HDI'': ''1''}'
$ws.Range("A133").Value = '{''RaceEthnicity'': ''1'', ''Hobby'': ''1'', ''SexualOrientation'': ''1''}'
$ws.Range("A134").Value = '{''Dependents'': ''2'', ''Hobby'': ''1'', ''SexualOrientation'': ''1''}'
$ws.Range("A135").Value = '{''Continent'': ''1'', ''Hobby'': ''1'', ''SexualOrientation'': ''1''}'
$ws.Range("A136").Value = '{''Age'': ''3'', ''Hobby'': ''1'', ''SexualOrientation'': ''1''}'
$ws.Range("A137").Value = '{''Hobby'': ''1'', ''SexualOrientation'': ''1'', ''HDI'': ''1''}'
$ws.Range("A138").Value = '{''Hobby'': ''1'', ''SexualOrientation'': ''1'', ''GDP'': ''1''}'
$ws.Range("A139").Value = '{''GINI'': ''1'', ''Hobby'': ''1'', ''SexualOrientation'': ''1''}'
$ws.Range("A140").Value = '{''GINI'': ''2'', ''Hobby'': ''1'', ''SexualOrientation'': ''1''}'
$ws.Range("A141").Value = '{''Dependents'': ''2'', ''Hobby'': ''1'', ''RaceEthnicity'': ''1''}'
$ws.Range("A143").Value = '{''Age'': ''3'', ''RaceEthnicity'': ''1'', ''Hobby'': ''1''}'
$ws.Range("A145").Value = '{''RaceEthnicity'': ''1'', ''Hobby'': ''1'', ''GDP'': ''1''}'
$ws.Range("A146").Value = '{''GINI'': ''2'', ''RaceEthnicity'': ''1'', ''Hobby'': ''1''}'
$ws.Range("A147").Value = '{''Age'': ''3'', ''Dependents'': ''2'', ''Hobby'': ''1''}'
$ws.Range("A148").Value = '{''Dependents'': ''2'', ''Hobby'': ''1'', ''HDI'': ''1''}'
$ws.Range("A149").Value = '{''GINI'': ''2'', ''Dependents'': ''2'', ''Hobby'': ''1''}'
$ws.Range("A151").Value = '{''Hobby'': ''1'', ''Continent'': ''1'', ''GDP'': ''1''}'
$ws.Range("A152").Value = '{''Age'': ''3'', ''Hobby'': ''1'', ''HDI'': ''1''}'
$ws.Range("A153").Value = '{''Hobby'': ''1'', ''HDI'': ''1'', ''GDP'': ''1''}'
$ws.Range("A155").Value = '{''GINI'': ''1'', ''Hobby'': ''1'', ''GDP'': ''1''}'
$ws.Range("A156").Value = '{''Student'': ''1'', ''Continent'': ''2'', ''Country'': ''2''}'
$ws.Range("A157").Value = '{''HDI'': ''1'', ''Country'': ''2'', ''Student'': ''1''}'
$ws.Range("A158").Value = '{''Student'': ''1'', ''Country'': ''2'', ''GDP'': ''2''}'
$ws.Range("A159").Value = '{''GINI'': ''2'', ''Student'': ''1'', ''Country'': ''2''}'
$ws.Range("A161").Value = '{''HDI'': ''1'', ''Country'': ''2'', ''SexualOrientation'': ''1''}'
$ws.Range("A162").Value = '{''Country'': ''2'', ''SexualOrientation'': ''1'', ''GDP'': ''2''}'
$ws.Range("A163").Value = '{''GINI'': ''2'', ''Country'': ''2'', ''SexualOrientation'': ''1''}'
$ws.Range("A164").Value = '{''HDI'': ''1'', ''Continent'': ''2'', ''Country'': ''2''}'
$ws.Range("A165").Value = '{''GDP'': ''2'', ''Continent'': ''2'', ''Country'': ''2''}'
$ws.Range("A166").Value = '{''GINI'': ''2'', ''Continent'': ''2'', ''Country'': ''2''}'
$ws.Range("A167").Value = '{''HDI'': ''1'', ''Country'': ''2'', ''GDP'': ''2''}'
$ws.Range("A169").Value = '{''GINI'': ''2'', ''Country'': ''2'', ''GDP'': ''2''}'
$ws.Range("A170").Value = '{''FormalEducation'': ''1'', ''Student'': ''1'', ''SexualOrientation'': ''1''}'
$ws.Range("A171").Value = '{''FormalEducation'': ''1'', ''RaceEthnicity'': ''1'', ''Student'': ''1''}'
$ws.Range("A172").Value = '{''FormalEducation'': ''1'', ''HDI'': ''1'', ''Student'': ''1''}'
$ws.Range("A173").Value = '{''UndergradMajor'': ''2'', ''Student'': ''1'', ''DevType'': ''2''}'
$ws.Range("A174").Value = '{''UndergradMajor'': ''2'', ''Student'': ''1'', ''SexualOrientation'': ''1''}'
$ws.Range("A176").Value = '{''UndergradMajor'': ''2'', ''Dependents'': ''2'', ''Student'': ''1''}'
$ws.Range("A177").Value = '{''UndergradMajor'': ''2'', ''HDI'': ''1'', ''Student'': ''1''}'
$ws.Range("A178").Value = '{''Student'': ''1'', ''DevType'': ''2'', ''SexualOrientation'': ''1''}'
$ws.Range("A179").Value = '{''RaceEthnicity'': ''1'', ''Student'': ''1'', ''DevType'': ''2''}'
$ws.Range("A180").Value = '{''Dependents'': ''2'', ''Student'': ''1'', ''DevType'': ''2''}'
$ws.Range("A181").Value = '{''HDI'': ''1'', ''Student'': ''1'', ''DevType'': ''2''}'
$ws.Range("A182").Value = '{''Exercise'': ''3'', ''Student'': ''1'', ''SexualOrientation'': ''1''}'
$ws.Range("A184").Value = '{''Dependents'': ''2'', ''Student'': ''1'', ''SexualOrientation'': ''1''}'
$ws.Range("A185").Value = '{''Continent'': ''1'', ''Student'': ''1'', ''SexualOrientation'': ''1''}'
$ws.Range("A186").Value = '{''Age'': ''3'', ''Student'': ''1'', ''SexualOrientation'': ''1''}'
$ws.Range("A187").Value = '{''HDI'': ''1'', ''Student'': ''1'', ''SexualOrientation'': ''1''}'
$ws.Range("A189").Value = '{''GINI'': ''1'', ''Student'': ''1'', ''SexualOrientation'': ''1''}'
$ws.Range("A190").Value = '{''GINI'': ''2'', ''Student'': ''1'', ''SexualOrientation'': ''1''}'
$ws.Range("A191").Value = '{''RaceEthnicity'': ''1'', ''Dependents'': ''2'', ''Student'': ''1''}'
$ws.Range("A193").Value = '{''Age'': ''3'', ''RaceEthnicity'': ''1'', ''Student'': ''1''}'
$ws.Range("A194").Value = '{''RaceEthnicity'': ''1'', ''HDI'': ''1'', ''Student'': ''1''}'
$ws.Range("A195").Value = '{''RaceEthnicity'': ''1'', ''GDP'': ''1'', ''Student'': ''1''}'
$ws.Range("A196").Value = '{''GINI'': ''2'', ''RaceEthnicity'': ''1'', ''Student'': ''1''}'
$ws.Range("A197").Value = '{''Age'': ''3'', ''Dependents'': ''2'', ''Student'': ''1''}'
$ws.Range("A198").Value = '{''Dependents'': ''2'', ''Student'': ''1'', ''HDI'': ''1''}'
$ws.Range("A199").Value = '{''HDI'': ''1'', ''Student'': ''1'', ''Continent'': ''1''}'
$ws.Range("A200").Value = '{''GDP'': ''1'', ''Student'': ''1'', ''Continent'': ''1''}'
$ws.Range("A201").Value = '{''Student'': ''1'', ''Continent'': ''2'', ''GDP'': ''2''}'
$ws.Range("A202").Value = '{''GINI'': ''2'', ''Student'': ''1'', ''Continent'': ''2''}'
$ws.Range("A203").Value = '{''Age'': ''3'', ''HDI'': ''1'', ''Student'': ''1''}'
$ws.Range("A204").Value = '{''GDP'': ''1'', ''Student'': ''1'', ''HDI'': ''1''}'
$ws.Range("A205").Value = '{''HDI'': ''1'', ''Student'': ''1'', ''GDP'': ''2''}'
$ws.Range("A206").Value = '{''GINI'': ''1'', ''HDI'': ''1'', ''Student'': ''1''}'
$ws.Range("A207").Value = '{''GINI'': ''2'', ''HDI'': ''1'', ''Student'': ''1''}'
$ws.Range("A208").Value = '{''GINI'': ''1'', ''GDP'': ''1'', ''Student'': ''1''}'
$ws.Range("A209").Value = '{''GINI'': ''2'', ''Student'': ''1'', ''GDP'': ''2''}'
$ws.Range("A210").Value = '{''FormalEducation'': ''1'', ''UndergradMajor'': ''2'', ''SexualOrientation'': ''1''}'
$ws.Range("A211").Value = '{''FormalEducation'': ''1'', ''RaceEthnicity'': ''1'', ''SexualOrientation'': ''1''}'
$ws.Range("A212").Value = '{''FormalEducation'': ''1'', ''Dependents'': ''2'', ''SexualOrientation'': ''1''}'
$ws.Range("A213").Value = '{''FormalEducation'': ''1'', ''HDI'': ''1'', ''SexualOrientation'': ''1''}'
$ws.Range("A214").Value = '{''GINI'': ''2'', ''FormalEducation'': ''1'', ''SexualOrientation'': ''1''}'
$ws.Range("A215").Value = '{''FormalEducation'': ''1'', ''HDI'': ''1'', ''RaceEthnicity'': ''1''}'
$ws.Range("A216").Value = '{''UndergradMajor'': ''2'', ''DevType'': ''2'', ''SexualOrientation'': ''1''}'
$ws.Range("A217").Value = '{''UndergradMajor'': ''2'', ''RaceEthnicity'': ''1'', ''DevType'': ''2''}'
$ws.Range("A218").Value = '{''UndergradMajor'': ''2'', ''HDI'': ''1'', ''DevType'': ''2''}'
$ws.Range("A219").Value = '{''UndergradMajor'': ''2'', ''RaceEthnicity'': ''1'', ''SexualOrientation'': ''1''}'
$ws.Range("A220").Value = '{''UndergradMajor'': ''2'', ''Dependents'': ''2'', ''SexualOrientation'': ''1''}'
$ws.Range("A221").Value = '{''Age'': ''3'', ''UndergradMajor'': ''2'', ''SexualOrientation'': ''1''}'
$ws.Range("A222").Value = '{''UndergradMajor'': ''2'', ''HDI'': ''1'', ''SexualOrientation'': ''1''}'
$ws.Range("A223").Value = '{''GINI'': ''2'', ''UndergradMajor'': ''2'', ''SexualOrientation'': ''1''}'
$ws.Range("A224").Value = '{''UndergradMajor'': ''2'', ''Dependents'': ''2'', ''RaceEthnicity'': ''1''}'
$ws.Range("A225").Value = '{''UndergradMajor'': ''2'', ''HDI'': ''1'', ''RaceEthnicity'': ''1''}'
$ws.Range("A227").Value = '{''UndergradMajor'': ''2'', ''HDI'': ''1'', ''GDP'': ''1''}'
$ws.Range("A229").Value = '{''Dependents'': ''2'', ''DevType'': ''2'', ''SexualOrientation'': ''1''}'
$ws.Range("A230").Value = '{''Age'': ''3'', ''DevType'': ''2'', ''SexualOrientation'': ''1''}'
$ws.Range("A231").Value = '{''HDI'': ''1'', ''DevType'': ''2'', ''SexualOrientation'': ''1''}'
$ws.Range("A232").Value = '{''GINI'': ''2'', ''DevType'': ''2'', ''SexualOrientation'': ''1''}'
$ws.Range("A233").Value = '{''RaceEthnicity'': ''1'', ''Dependents'': ''2'', ''DevType'': ''2''}'
$ws.Range("A234").Value = '{''RaceEthnicity'': ''1'', ''HDI'': ''1'', ''DevType'': ''2''}'
$ws.Range("A235").Value = '{''Dependents'': ''2'', ''HDI'': ''1'', ''DevType'': ''2''}'
$ws.Range("A236").Value = '{''RaceEthnicity'': ''1'', ''Dependents'': ''2'', ''SexualOrientation'': ''1''}'
$ws.Range("A237").Value = '{''Continent'': ''1'', ''RaceEthnicity'': ''1'', ''SexualOrientation'': ''1''}'
$ws.Range("A238").Value = '{''Age'': ''3'', ''RaceEthnicity'': ''1'', ''SexualOrientation'': ''1''}'
$ws.Range("A239").Value = '{''RaceEthnicity'': ''1'', ''HDI'': ''1'', ''SexualOrientation'': ''1''}'
$ws.Range("A241").Value = '{''GINI'': ''1'', ''RaceEthnicity'': ''1'', ''SexualOrientation'': ''1''}'
$ws.Range("A242").Value = '{''GINI'': ''2'', ''RaceEthnicity'': ''1'', ''SexualOrientation'': ''1''}'
$ws.Range("A243").Value = '{''Continent'': ''1'', ''Dependents'': ''2'', ''SexualOrientation'': ''1''}'
$ws.Range("A244").Value = '{''Age'': ''3'', ''Dependents'': ''2'', ''SexualOrientation'': ''1''}'
$ws.Range("A245").Value = '{''Dependents'': ''2'', ''HDI'': ''1'', ''SexualOrientation'': ''1''}'
$ws.Range("A246").Value = '{''Dependents'': ''2'', ''GDP'': ''1'', ''SexualOrientation'': ''1''}'
$ws.Range("A247").Value = '{''GINI'': ''2'', ''Dependents'': ''2'', ''SexualOrientation'': ''1''}'
$ws.Range("A248").Value = '{''Continent'': ''1'', ''HDI'': ''1'', ''SexualOrientation'': ''1''}'
$ws.Range("A249").Value = '{''Continent'': ''1'', ''GDP'': ''1'', ''SexualOrientation'': ''1''}'
$ws.Range("A251").Value = '{''GINI'': ''2'', ''Continent'': ''2'', ''SexualOrientation'': ''1''}'
$ws.Range("A252").Value = '{''Age'': ''3'', ''HDI'': ''1'', ''SexualOrientation'': ''1''}'
$ws.Range("A253").Value = '{''HDI'': ''1'', ''SexualOrientation'': ''1'', ''GDP'': ''1''}'
$ws.Range("A254").Value = '{''HDI'': ''1'', ''SexualOrientation'': ''1'', ''GDP'': ''2''}'
$ws.Range("A255").Value = '{''GINI'': ''1'', ''HDI'': ''1'', ''SexualOrientation'': ''1''}'
$ws.Range("A256").Value = '{''GINI'': ''1'', ''GDP'': ''1'', ''SexualOrientation'': ''1''}'
$ws.Range("A257").Value = '{''GINI'': ''2'', ''SexualOrientation'': ''1'', ''GDP'': ''2''}'
$ws.Range("A258").Value = '{''Age'': ''3'', ''RaceEthnicity'': ''1'', ''Dependents'': ''2''}'
$ws.Range("A259").Value = '{''RaceEthnicity'': ''1'', ''Dependents'': ''2'', ''HDI'': ''1''}'
$ws.Range("A260").Value = '{''RaceEthnicity'': ''1'', ''HDI'': ''1'', ''Continent'': ''1''}'
$ws.Range("A262").Value = '{''Age'': ''3'', ''RaceEthnicity'': ''1'', ''HDI'': ''1''}'
$ws.Range("A264").Value = '{''GINI'': ''1'', ''RaceEthnicity'': ''1'', ''HDI'': ''1''}'
$ws.Range("A265").Value = '{''GINI'': ''1'', ''RaceEthnicity'': ''1'', ''GDP'': ''1''}'
$ws.Range("A266").Value = '{''Age'': ''3'', ''Dependents'': ''2'', ''HDI'': ''1''}'
$ws.Range("A267").Value = '{''Dependents'': ''2'', ''GDP'': ''1'', ''HDI'': ''1''}'
$ws.Range("A268").Value = '{''HDI'': ''1'', ''Continent'': ''1'', ''GDP'': ''1''}'
$ws.Range("A269").Value = '{''HDI'': ''1'', ''Continent'': ''2'', ''GDP'': ''2''}'
$ws.Range("A270").Value = '{''GINI'': ''2'', ''HDI'': ''1'', ''Continent'': ''2''}'
$ws.Range("A271").Value = '{''GINI'': ''2'', ''Continent'': ''2'', ''GDP'': ''2''}'
$ws.Range("A272").Value = '{''GINI'': ''1'', ''HDI'': ''1'', ''GDP'': ''1''}'
$ws.Range("A273").Value = '{''GINI'': ''2'', ''HDI'': ''1'', ''GDP'': ''2''}'
$ws.Range("A274").Value = '{''FormalEducation'': ''1'', ''Hobby'': ''1'', ''Student'': ''1'', ''SexualOrientation'': ''1''}'
$ws.Range("A275").Value = '{''UndergradMajor'': ''2'', ''Hobby'': ''1'', ''Student'': ''1'', ''SexualOrientation'': ''1''}'
$ws.Range("A276").Value = '{''RaceEthnicity'': ''1'', ''Hobby'': ''1'', ''Student'': ''1'', ''UndergradMajor'': ''2''}'
$ws.Range("A277").Value = '{''Dependents'': ''2'', ''Hobby'': ''1'', ''Student'': ''1'', ''UndergradMajor'': ''2''}'
$ws.Range("A278").Value = '{''UndergradMajor'': ''2'', ''Hobby'': ''1'', ''Student'': ''1'', ''HDI'': ''1''}'
$ws.Range("A279").Value = '{''Hobby'': ''1'', ''Student'': ''1'', ''DevType'': ''2'', ''SexualOrientation'': ''1''}'
$ws.Range("A280").Value = '{''RaceEthnicity'': ''1'', ''Hobby'': ''1'', ''Student'': ''1'', ''DevType'': ''2''}'
$ws.Range("A281").Value = '{''Hobby'': ''1'', ''Student'': ''1'', ''DevType'': ''2'', ''HDI'': ''1''}'
$ws.Range("A282").Value = '{''RaceEthnicity'': ''1'', ''Hobby'': ''1'', ''Student'': ''1'', ''SexualOrientation'': ''1''}'
$ws.Range("A283").Value = '{''Dependents'': ''2'', ''Hobby'': ''1'', ''Student'': ''1'', ''SexualOrientation'': ''1''}'
$ws.Range("A284").Value = '{''Age'': ''3'', ''Hobby'': ''1'', ''Student'': ''1'', ''SexualOrientation'': ''1''}'
$ws.Range("A285").Value = '{''Hobby'': ''1'', ''Student'': ''1'', ''SexualOrientation'': ''1'', ''HDI'': ''1''}'
$ws.Range("A286").Value = '{''Hobby'': ''1'', ''Student'': ''1'', ''SexualOrientation'': ''1'', ''GDP'': ''1''}'
$ws.Range("A287").Value = '{''GINI'': ''2'', ''Hobby'': ''1'', ''Student'': ''1'', ''SexualOrientation'': ''1''}'
$ws.Range("A288").Value = '{''RaceEthnicity'': ''1'', ''Hobby'': ''1'', ''Student'': ''1'', ''Dependents'': ''2''}'
$ws.Range("A289").Value = '{''Age'': ''3'', ''RaceEthnicity'': ''1'', ''Hobby'': ''1'', ''Student'': ''1''}'
$ws.Range("A290").Value = '{''RaceEthnicity'': ''1'', ''Hobby'': ''1'', ''Student'': ''1'', ''HDI'': ''1''}'
$ws.Range("A291").Value = '{''Dependents'': ''2'', ''Hobby'': ''1'', ''Student'': ''1'', ''HDI'': ''1''}'
$ws.Range("A292").Value = '{''Age'': ''3'', ''Hobby'': ''1'', ''Student'': ''1'', ''HDI'': ''1''}'
$ws.Range("A293").Value = '{''Hobby'': ''1'', ''Student'': ''1'', ''HDI'': ''1'', ''GDP'': ''1''}'
$ws.Range("A294").Value = '{''UndergradMajor'': ''2'', ''Hobby'': ''1'', ''DevType'': ''2'', ''SexualOrientation'': ''1''}'
$ws.Range("A295").Value = '{''UndergradMajor'': ''2'', ''Hobby'': ''1'', ''RaceEthnicity'': ''1'', ''SexualOrientation'': ''1''}'
$ws.Range("A296").Value = '{''Dependents'': ''2'', ''UndergradMajor'': ''2'', ''Hobby'': ''1'', ''SexualOrientation'': ''1''}'
$ws.Range("A297").Value = '{''UndergradMajor'': ''2'', ''Hobby'': ''1'', ''SexualOrientation'': ''1'', ''HDI'': ''1''}'
$ws.Range("A298").Value = '{''Dependents'': ''2'', ''Hobby'': ''1'', ''RaceEthnicity'': ''1'', ''UndergradMajor'': ''2''}'
$ws.Range("A299").Value = '{''UndergradMajor'': ''2'', ''Hobby'': ''1'', ''RaceEthnicity'': ''1'', ''HDI'': ''1''}'
$ws.Range("A300").Value = '{''Dependents'': ''2'', ''UndergradMajor'': ''2'', ''Hobby'': ''1'', ''HDI'': ''1''}'
$ws.Range("A301").Value = '{''RaceEthnicity'': ''1'', ''Hobby'': ''1'', ''DevType'': ''2'', ''SexualOrientation'': ''1''}'
$ws.Range("A302").Value = '{''Dependents'': ''2'', ''Hobby'': ''1'', ''DevType'': ''2'', ''SexualOrientation'': ''1''}'
$ws.Range("A303").Value = '{''Hobby'': ''1'', ''DevType'': ''2'', ''SexualOrientation'': ''1'', ''HDI'': ''1''}'
$ws.Range("A304").Value = '{''RaceEthnicity'': ''1'', ''Hobby'': ''1'', ''DevType'': ''2'', ''HDI'': ''1''}'
$ws.Range("A305").Value = '{''Dependents'': ''2'', ''Hobby'': ''1'', ''RaceEthnicity'': ''1'', ''SexualOrientation'': ''1''}'
$ws.Range("A306").Value = '{''Continent'': ''1'', ''RaceEthnicity'': ''1'', ''Hobby'': ''1'', ''SexualOrientation'': ''1''}'
$ws.Range("A307").Value = '{''Age'': ''3'', ''RaceEthnicity'': ''1'', ''Hobby'': ''1'', ''SexualOrientation'': ''1''}'
$ws.Range("A308").Value = '{''RaceEthnicity'': ''1'', ''Hobby'': ''1'', ''SexualOrientation'': ''1'', ''HDI'': ''1''}'
$ws.Range("A309").Value = '{''RaceEthnicity'': ''1'', ''Hobby'': ''1'', ''SexualOrientation'': ''1'', ''GDP'': ''1''}'
$ws.Range("A310").Value = '{''GINI'': ''2'', ''RaceEthnicity'': ''1'', ''Hobby'': ''1'', ''SexualOrientation'': ''1''}'
$ws.Range("A311").Value = '{''Age'': ''3'', ''Dependents'': ''2'', ''Hobby'': ''1'', ''SexualOrientation'': ''1''}'
$ws.Range("A312").Value = '{''Dependents'': ''2'', ''Hobby'': ''1'', ''SexualOrientation'': ''1'', ''HDI'': ''1''}'
$ws.Range("A313").Value = '{''Age'': ''3'', ''Hobby'': ''1'', ''SexualOrientation'': ''1'', ''HDI'': ''1''}'
$ws.Range("A314").Value = '{''Hobby'': ''1'', ''GDP'': ''1'', ''SexualOrientation'': ''1'', ''HDI'': ''1''}'
$ws.Range("A315").Value = '{''Dependents'': ''2'', ''Hobby'': ''1'', ''RaceEthnicity'': ''1'', ''HDI'': ''1''}'
$ws.Range("A316").Value = '{''RaceEthnicity'': ''1'', ''Hobby'': ''1'', ''HDI'': ''1'', ''GDP'': ''1''}'
$ws.Range("A317").Value = '{''Continent'': ''1'', ''Hobby'': ''1'', ''HDI'': ''1'', ''GDP'': ''1''}'
$ws.Range("A318").Value = '{''GINI'': ''1'', ''Hobby'': ''1'', ''HDI'': ''1'', ''GDP'': ''1''}'
$ws.Range("A319").Value = '{''Student'': ''1'', ''HDI'': ''1'', ''Continent'': ''2'', ''Country'': ''2''}'
$ws.Range("A320").Value = '{''Student'': ''1'', ''Continent'': ''2'', ''GDP'': ''2'', ''Country'': ''2''}'
$ws.Range("A321").Value = '{''GINI'': ''2'', ''Student'': ''1'', ''Continent'': ''2'', ''Country'': ''2''}'
$ws.Range("A322").Value = '{''HDI'': ''1'', ''Country'': ''2'', ''GDP'': ''2'', ''Student'': ''1''}'
$ws.Range("A323").Value = '{''GINI'': ''2'', ''HDI'': ''1'', ''Country'': ''2'', ''Student'': ''1''}'
$ws.Range("A324").Value = '{''GINI'': ''2'', ''Student'': ''1'', ''Country'': ''2'', ''GDP'': ''2''}'
$ws.Range("A325").Value = '{''HDI'': ''1'', ''Continent'': ''2'', ''SexualOrientation'': ''1'', ''Country'': ''2''}'
$ws.Range("A326").Value = '{''GDP'': ''2'', ''Continent'': ''2'', ''SexualOrientation'': ''1'', ''Country'': ''2''}'
$ws.Range("A327").Value = '{''GINI'': ''2'', ''Continent'': ''2'', ''SexualOrientation'': ''1'', ''Country'': ''2''}'
$ws.Range("A328").Value = '{''HDI'': ''1'', ''Country'': ''2'', ''SexualOrientation'': ''1'', ''GDP'': ''2''}'
$ws.Range("A329").Value = '{''GINI'': ''2'', ''HDI'': ''1'', ''Country'': ''2'', ''SexualOrientation'': ''1''}'
$ws.Range("A330").Value = '{''GINI'': ''2'', ''Country'': ''2'', ''SexualOrientation'': ''1'', ''GDP'': ''2''}'
$ws.Range("A331").Value = '{''HDI'': ''1'', ''Continent'': ''2'', ''GDP'': ''2'', ''Country'': ''2''}'
$ws.Range("A332").Value = '{''GINI'': ''2'', ''HDI'': ''1'', ''Continent'': ''2'', ''Country'': ''2''}'
$ws.Range("A333").Value = '{''GINI'': ''2'', ''GDP'': ''2'', ''Continent'': ''2'', ''Country'': ''2''}'
$ws.Range("A334").Value = '{''GINI'': ''2'', ''HDI'': ''1'', ''Country'': ''2'', ''GDP'': ''2''}'
$ws.Range("A335").Value = '{''FormalEducation'': ''1'', ''RaceEthnicity'': ''1'', ''Student'': ''1'', ''SexualOrientation'': ''1''}'
$ws.Range("A336").Value = '{''FormalEducation'': ''1'', ''HDI'': ''1'', ''Student'': ''1'', ''SexualOrientation'': ''1''}'
$ws.Range("A337").Value = '{''UndergradMajor'': ''2'', ''Student'': ''1'', ''DevType'': ''2'', ''SexualOrientation'': ''1''}'
$ws.Range("A338").Value = '{''UndergradMajor'': ''2'', ''RaceEthnicity'': ''1'', ''Student'': ''1'', ''SexualOrientation'': ''1''}'
$ws.Range("A339").Value = '{''UndergradMajor'': ''2'', ''Dependents'': ''2'', ''Student'': ''1'', ''SexualOrientation'': ''1''}'
$ws.Range("A340").Value = '{''UndergradMajor'': ''2'', ''HDI'': ''1'', ''Student'': ''1'', ''SexualOrientation'': ''1''}'
$ws.Range("A341").Value = '{''RaceEthnicity'': ''1'', ''HDI'': ''1'', ''Student'': ''1'', ''UndergradMajor'': ''2''}'
$ws.Range("A342").Value = '{''RaceEthnicity'': ''1'', ''Student'': ''1'', ''DevType'': ''2'', ''SexualOrientation'': ''1''}'
$ws.Range("A343").Value = '{''Dependents'': ''2'', ''Student'': ''1'', ''DevType'': ''2'', ''SexualOrientation'': ''1''}'
$ws.Range("A344").Value = '{''HDI'': ''1'', ''Student'': ''1'', ''DevType'': ''2'', ''SexualOrientation'': ''1''}'
$ws.Range("A345").Value = '{''RaceEthnicity'': ''1'', ''HDI'': ''1'', ''Student'': ''1'', ''DevType'': ''2''}'
$ws.Range("A346").Value = '{''RaceEthnicity'': ''1'', ''Dependents'': ''2'', ''Student'': ''1'', ''SexualOrientation'': ''1''}'
$ws.Range("A347").Value = '{''Continent'': ''1'', ''RaceEthnicity'': ''1'', ''Student'': ''1'', ''SexualOrientation'': ''1''}'
$ws.Range("A348").Value = '{''Age'': ''3'', ''RaceEthnicity'': ''1'', ''Student'': ''1'', ''SexualOrientation'': ''1''}'
$ws.Range("A349").Value = '{''RaceEthnicity'': ''1'', ''HDI'': ''1'', ''Student'': ''1'', ''SexualOrientation'': ''1''}'
$ws.Range("A351").Value = '{''GINI'': ''2'', ''RaceEthnicity'': ''1'', ''Student'': ''1'', ''SexualOrientation'': ''1''}'
$ws.Range("A352").Value = '{''Age'': ''3'', ''Dependents'': ''2'', ''Student'': ''1'', ''SexualOrientation'': ''1''}'
$ws.Range("A353").Value = '{''Dependents'': ''2'', ''Student'': ''1'', ''SexualOrientation'': ''1'', ''HDI'': ''1''}'
$ws.Range("A354").Value = '{''Age'': ''3'', ''HDI'': ''1'', ''Student'': ''1'', ''SexualOrientation'': ''1''}'
$ws.Range("A356").Value = '{''RaceEthnicity'': ''1'', ''Dependents'': ''2'', ''Student'': ''1'', ''HDI'': ''1''}'
$ws.Range("A357").Value = '{''Age'': ''3'', ''RaceEthnicity'': ''1'', ''HDI'': ''1'', ''Student'': ''1''}'
$ws.Range("A358").Value = '{''RaceEthnicity'': ''1'', ''GDP'': ''1'', ''Student'': ''1'', ''HDI'': ''1''}'
$ws.Range("A359").Value = '{''Continent'': ''1'', ''GDP'': ''1'', ''Student'': ''1'', ''HDI'': ''1''}'
$ws.Range("A360").Value = '{''HDI'': ''1'', ''Continent'': ''2'', ''GDP'': ''2'', ''Student'': ''1''}'
$ws.Range("A361").Value = '{''GINI'': ''2'', ''HDI'': ''1'', ''Continent'': ''2'', ''Student'': ''1''}'
$ws.Range("A362").Value = '{''GINI'': ''2'', ''Student'': ''1'', ''Continent'': ''2'', ''GDP'': ''2''}'
$ws.Range("A363").Value = '{''GINI'': ''1'', ''GDP'': ''1'', ''Student'': ''1'', ''HDI'': ''1''}'
$ws.Range("A364").Value = '{''GINI'': ''2'', ''HDI'': ''1'', ''Student'': ''1'', ''GDP'': ''2''}'
$ws.Range("A365").Value = '{''FormalEducation'': ''1'', ''HDI'': ''1'', ''RaceEthnicity'': ''1'', ''SexualOrientation'': ''1''}'
$ws.Range("A366").Value = '{''UndergradMajor'': ''2'', ''RaceEthnicity'': ''1'', ''DevType'': ''2'', ''SexualOrientation'': ''1''}'
$ws.Range("A367").Value = '{''UndergradMajor'': ''2'', ''HDI'': ''1'', ''DevType'': ''2'', ''SexualOrientation'': ''1''}'
$ws.Range("A368").Value = '{''UndergradMajor'': ''2'', ''Dependents'': ''2'', ''RaceEthnicity'': ''1'', ''SexualOrientation'': ''1''}'
$ws.Range("A369").Value = '{''UndergradMajor'': ''2'', ''HDI'': ''1'', ''RaceEthnicity'': ''1'', ''SexualOrientation'': ''1''}'
$ws.Range("A370").Value = '{''Dependents'': ''2'', ''UndergradMajor'': ''2'', ''HDI'': ''1'', ''SexualOrientation'': ''1''}'
$ws.Range("A371").Value = '{''Dependents'': ''2'', ''HDI'': ''1'', ''RaceEthnicity'': ''1'', ''UndergradMajor'': ''2''}'
$ws.Range("A372").Value = '{''RaceEthnicity'': ''1'', ''Dependents'': ''2'', ''DevType'': ''2'', ''SexualOrientation'': ''1''}'
$ws.Range("A373").Value = '{''RaceEthnicity'': ''1'', ''HDI'': ''1'', ''DevType'': ''2'', ''SexualOrientation'': ''1''}'
$ws.Range("A374").Value = '{''Dependents'': ''2'', ''HDI'': ''1'', ''DevType'': ''2'', ''SexualOrientation'': ''1''}'
$ws.Range("A375").Value = '{''RaceEthnicity'': ''1'', ''Dependents'': ''2'', ''SexualOrientation'': ''1'', ''HDI'': ''1''}'
$ws.Range("A376").Value = '{''Continent'': ''1'', ''RaceEthnicity'': ''1'', ''HDI'': ''1'', ''SexualOrientation'': ''1''}'
$ws.Range("A377").Value = '{''Continent'': ''1'', ''RaceEthnicity'': ''1'', ''GDP'': ''1'', ''SexualOrientation'': ''1''}'
$ws.Range("A378").Value = '{''Age'': ''3'', ''RaceEthnicity'': ''1'', ''HDI'': ''1'', ''SexualOrientation'': ''1''}'
$ws.Range("A380").Value = '{''GINI'': ''1'', ''RaceEthnicity'': ''1'', ''HDI'': ''1'', ''SexualOrientation'': ''1''}'
$ws.Range("A381").Value = '{''GINI'': ''1'', ''RaceEthnicity'': ''1'', ''GDP'': ''1'', ''SexualOrientation'': ''1''}'
$ws.Range("A382").Value = '{''Dependents'': ''2'', ''GDP'': ''1'', ''SexualOrientation'': ''1'', ''HDI'': ''1''}'
$ws.Range("A383").Value = '{''Continent'': ''1'', ''HDI'': ''1'', ''SexualOrientation'': ''1'', ''GDP'': ''1''}'
$ws.Range("A384").Value = '{''HDI'': ''1'', ''Continent'': ''2'', ''SexualOrientation'': ''1'', ''GDP'': ''2''}'
$ws.Range("A385").Value = '{''GINI'': ''2'', ''HDI'': ''1'', ''Continent'': ''2'', ''SexualOrientation'': ''1''}'
$ws.Range("A386").Value = '{''GINI'': ''2'', ''Continent'': ''2'', ''SexualOrientation'': ''1'', ''GDP'': ''2''}'
$ws.Range("A387").Value = '{''GINI'': ''1'', ''HDI'': ''1'', ''SexualOrientation'': ''1'', ''GDP'': ''1''}'
$ws.Range("A388").Value = '{''GINI'': ''2'', ''HDI'': ''1'', ''SexualOrientation'': ''1'', ''GDP'': ''2''}'
$ws.Range("A389").Value = '{''Continent'': ''1'', ''RaceEthnicity'': ''1'', ''GDP'': ''1'', ''HDI'': ''1''}'
$ws.Range("A390").Value = '{''GINI'': ''1'', ''RaceEthnicity'': ''1'', ''GDP'': ''1'', ''HDI'': ''1''}'
$ws.Range("A391").Value = '{''GINI'': ''2'', ''HDI'': ''1'', ''Continent'': ''2'', ''GDP'': ''2''}'
$ws.Range("A392").Value = '{''Hobby'': ''1'', ''SexualOrientation'': ''1'', ''UndergradMajor'': ''2'', ''RaceEthnicity'': ''1'', ''Student'': ''1''}'
$ws.Range("A393").Value = '{''Hobby'': ''1'', ''SexualOrientation'': ''1'', ''UndergradMajor'': ''2'', ''HDI'': ''1'', ''Student'': ''1''}'
$ws.Range("A394").Value = '{''Hobby'': ''1'', ''UndergradMajor'': ''2'', ''RaceEthnicity'': ''1'', ''HDI'': ''1'', ''Student'': ''1''}'
$ws.Range("A395").Value = '{''Hobby'': ''1'', ''DevType'': ''2'', ''SexualOrientation'': ''1'', ''RaceEthnicity'': ''1'', ''Student'': ''1''}'
$ws.Range("A396").Value = '{''Hobby'': ''1'', ''DevType'': ''2'', ''SexualOrientation'': ''1'', ''HDI'': ''1'', ''Student'': ''1''}'
$ws.Range("A397").Value = '{''Hobby'': ''1'', ''DevType'': ''2'', ''RaceEthnicity'': ''1'', ''HDI'': ''1'', ''Student'': ''1''}'
$ws.Range("A398").Value = '{''Hobby'': ''1'', ''SexualOrientation'': ''1'', ''RaceEthnicity'': ''1'', ''Dependents'': ''2'', ''Student'': ''1''}'
$ws.Range("A399").Value = '{''Hobby'': ''1'', ''SexualOrientation'': ''1'', ''RaceEthnicity'': ''1'', ''HDI'': ''1'', ''Student'': ''1''}'
$ws.Range("A400").Value = '{''Hobby'': ''1'', ''SexualOrientation'': ''1'', ''Dependents'': ''2'', ''Student'': ''1'', ''HDI'': ''1''}'
$ws.Range("A401").Value = '{''Hobby'': ''1'', ''SexualOrientation'': ''1'', ''HDI'': ''1'', ''Student'': ''1'', ''GDP'': ''1''}'
$ws.Range("A402").Value = '{''Hobby'': ''1'', ''RaceEthnicity'': ''1'', ''Dependents'': ''2'', ''Student'': ''1'', ''HDI'': ''1''}'
$ws.Range("A403").Value = '{''Hobby'': ''1'', ''SexualOrientation'': ''1'', ''UndergradMajor'': ''2'', ''RaceEthnicity'': ''1'', ''HDI'': ''1''}'
$ws.Range("A404").Value = '{''Hobby'': ''1'', ''DevType'': ''2'', ''SexualOrientation'': ''1'', ''RaceEthnicity'': ''1'', ''HDI'': ''1''}'
$ws.Range("A405").Value = '{''Hobby'': ''1'', ''SexualOrientation'': ''1'', ''RaceEthnicity'': ''1'', ''Dependents'': ''2'', ''HDI'': ''1''}'
$ws.Range("A406").Value = '{''Hobby'': ''1'', ''SexualOrientation'': ''1'', ''RaceEthnicity'': ''1'', ''HDI'': ''1'', ''GDP'': ''1''}'
$ws.Range("A407").Value = '{''Continent'': ''2'', ''HDI'': ''1'', ''Student'': ''1'', ''GDP'': ''2'', ''Country'': ''2''}'
$ws.Range("A408").Value = '{''Continent'': ''2'', ''HDI'': ''1'', ''Student'': ''1'', ''GINI'': ''2'', ''Country'': ''2''}'
$ws.Range("A409").Value = '{''Continent'': ''2'', ''Student'': ''1'', ''GDP'': ''2'', ''GINI'': ''2'', ''Country'': ''2''}'
$ws.Range("A410").Value = '{''HDI'': ''1'', ''Student'': ''1'', ''GDP'': ''2'', ''GINI'': ''2'', ''Country'': ''2''}'
$ws.Range("A411").Value = '{''SexualOrientation'': ''1'', ''Continent'': ''2'', ''HDI'': ''1'', ''GDP'': ''2'', ''Country'': ''2''}'
$ws.Range("A412").Value = '{''SexualOrientation'': ''1'', ''Continent'': ''2'', ''HDI'': ''1'', ''GINI'': ''2'', ''Country'': ''2''}'
$ws.Range("A413").Value = '{''SexualOrientation'': ''1'', ''Continent'': ''2'', ''GDP'': ''2'', ''GINI'': ''2'', ''Country'': ''2''}'
$ws.Range("A414").Value = '{''SexualOrientation'': ''1'', ''HDI'': ''1'', ''GDP'': ''2'', ''GINI'': ''2'', ''Country'': ''2''}'
$ws.Range("A415").Value = '{''Continent'': ''2'', ''HDI'': ''1'', ''GDP'': ''2'', ''GINI'': ''2'', ''Country'': ''2''}'
$ws.Range("A416").Value = '{''SexualOrientation'': ''1'', ''UndergradMajor'': ''2'', ''RaceEthnicity'': ''1'', ''HDI'': ''1'', ''Student'': ''1''}'
$ws.Range("A417").Value = '{''DevType'': ''2'', ''SexualOrientation'': ''1'', ''RaceEthnicity'': ''1'', ''HDI'': ''1'', ''Student'': ''1''}'
$ws.Range("A418").Value = '{''SexualOrientation'': ''1'', ''RaceEthnicity'': ''1'', ''Dependents'': ''2'', ''Student'': ''1'', ''HDI'': ''1''}'
$ws.Range("A419").Value = '{''SexualOrientation'': ''1'', ''RaceEthnicity'': ''1'', ''HDI'': ''1'', ''Student'': ''1'', ''GDP'': ''1''}'
$ws.Range("A420").Value = '{''Continent'': ''2'', ''HDI'': ''1'', ''Student'': ''1'', ''GDP'': ''2'', ''GINI'': ''2''}'
$ws.Range("A421").Value = '{''SexualOrientation'': ''1'', ''Continent'': ''1'', ''RaceEthnicity'': ''1'', ''HDI'': ''1'', ''GDP'': ''1''}'
$ws.Range("A422").Value = '{''SexualOrientation'': ''1'', ''GINI'': ''1'', ''RaceEthnicity'': ''1'', ''HDI'': ''1'', ''GDP'': ''1''}'
$ws.Range("A423").Value = '{''SexualOrientation'': ''1'', ''Continent'': ''2'', ''HDI'': ''1'', ''GDP'': ''2'', ''GINI'': ''2''}'
$ws.Range("A424").Value = '{''Hobby'': ''1'', ''SexualOrientation'': ''1'', ''RaceEthnicity'': ''1'', ''Dependents'': ''2'', ''Student'': ''1'', ''HDI'': ''1''}'
$ws.Range("A425").Value = '{''Continent'': ''2'', ''HDI'': ''1'', ''Student'': ''1'', ''GDP'': ''2'', ''GINI'': ''2'', ''Country'': ''2''}'
$ws.Range("A426").Value = '{''SexualOrientation'': ''1'', ''Continent'': ''2'', ''HDI'': ''1'', ''GDP'': ''2'', ''GINI'': ''2'', ''Country'': ''2''}'
